# Update NATMI LR-pair sheet with new TPM-derived values.
# The Sending cluster (column A) changes from "Resolving-Mac" to "ECs" for every
# data row, and the numeric expression/specificity columns (G, H, M, N, O, P, Q,
# R, S, T) are refreshed with newly recomputed TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("A2").Value = "ECs"
$ws.Range("G2").Value = 0.08586166666666667
$ws.Range("H2").Value = 0.257585
$ws.Range("M2").Value = 2.808848
$ws.Range("N2").Value = 8.426544
$ws.Range("O2").Value = 0.1240735410877844
$ws.Range("P2").Value = 0.1265382238946467
$ws.Range("Q2").Value = 0.2411723706933333
$ws.Range("R2").Value = 2.17055133624
$ws.Range("S2").Value = 0.1240735410877844
$ws.Range("T2").Value = 0.1265382238946467

# Row 3 (Target cluster: FAPs)
$ws.Range("A3").Value = "ECs"
$ws.Range("G3").Value = 0.08586166666666667
$ws.Range("H3").Value = 0.257585
$ws.Range("O3").Value = 0.3534602246216307
$ws.Range("P3").Value = 0.3604816034820782
$ws.Range("Q3").Value = 0.6870509181122223
$ws.Range("R3").Value = 6.18345826301
$ws.Range("S3").Value = 0.3534602246216307
$ws.Range("T3").Value = 0.3604816034820782

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("A4").Value = "ECs"
$ws.Range("G4").Value = 0.08586166666666667
$ws.Range("H4").Value = 0.257585
$ws.Range("M4").Value = 5.550351333333333
$ws.Range("N4").Value = 16.651054
$ws.Range("O4").Value = 0.2451723070126871
$ws.Range("P4").Value = 0.2500425796309677
$ws.Range("Q4").Value = 0.4765624160655555
$ws.Range("R4").Value = 4.28906174459
$ws.Range("S4").Value = 0.2451723070126871
$ws.Range("T4").Value = 0.2500425796309677

# Row 5 (Target cluster: MuSCs)
$ws.Range("A5").Value = "ECs"
$ws.Range("G5").Value = 0.08586166666666667
$ws.Range("H5").Value = 0.257585
$ws.Range("M5").Value = 1.322847
$ws.Range("N5").Value = 2.645694
$ws.Range("O5").Value = 0.0584333191427063
$ws.Range("P5").Value = 0.0397293860601361
$ws.Range("Q5").Value = 0.113581848165
$ws.Range("R5").Value = 0.68149108899
$ws.Range("S5").Value = 0.0584333191427063
$ws.Range("T5").Value = 0.0397293860601361

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("A6").Value = "ECs"
$ws.Range("G6").Value = 0.08586166666666667
$ws.Range("H6").Value = 0.257585
$ws.Range("M6").Value = 4.954692000000001
$ws.Range("N6").Value = 14.864076
$ws.Range("O6").Value = 0.2188606081351916
$ws.Range("P6").Value = 0.2232082069321711
$ws.Range("Q6").Value = 0.42541811294
$ws.Range("R6").Value = 3.82876301646
$ws.Range("S6").Value = 0.2188606081351916
$ws.Range("T6").Value = 0.2232082069321711
